# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right before the existing "2022-Q3"
#    sheet and fill it with the quarter's fund-holding detail rows.
# 2) Insert a new summary row at the top of the "总计" data (row 2) that
#    captures the 2022-Q4 totals, pushing the older quarters down and
#    renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q4" sheet before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $q4.Cells.Item(1, $col + 2).Value = $headers[$col]
}

# code, name, scale, position, ratio, marketValue, rank
$rows = @(
    @("008132", "鹏华价值驱动混合", "4.41", "93.64", "4.21", "0.1857", 5),
    @("014541", "华安新能源主题混合A", "1.09", "90.05", "3.05", "0.0332", 10),
    @("011888", "民生加银周期优选混合型证券投资基金A", "0.35", "92.72", "3.12", "0.0109", 10),
    @("006369", "弘毅远方国企转型升级混合A", "0.42", "77.78", "2.38", "0.0100", 7),
    @("620004", "金元顺安价值增长混合", "0.52", "77.68", "1.72", "0.0089", 8),
    @("014542", "华安新能源主题混合C", "0.10", "90.05", "3.05", "0.0030", 10),
    @("011889", "民生加银周期优选混合型证券投资基金C", "0.09", "92.72", "3.12", "0.0028", 10),
    @("013530", "弘毅远方国企转型升级混合C", "0.03", "77.78", "2.38", "0.0007", 7)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2

    $q4.Cells.Item($excelRow, 1).Value = $r

    # Columns B,D,E,F,G hold numeric-looking text (fund codes, percentages,
    # …) that must stay text (leading zeros / trailing zeros matter), so
    # force the cell to text format before writing the value.
    $q4.Cells.Item($excelRow, 2).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 2).Value = $rowData[0]

    $q4.Cells.Item($excelRow, 3).Value = $rowData[1]

    $q4.Cells.Item($excelRow, 4).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 4).Value = $rowData[2]

    $q4.Cells.Item($excelRow, 5).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 5).Value = $rowData[3]

    $q4.Cells.Item($excelRow, 6).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 6).Value = $rowData[4]

    $q4.Cells.Item($excelRow, 7).NumberFormat = "@"
    $q4.Cells.Item($excelRow, 7).Value = $rowData[5]

    $q4.Cells.Item($excelRow, 8).Value = $rowData[6]
}

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 summary row into "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 0.26

# Renumber the index column (A) for the quarters that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
